$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- new shared-string text blocks (single-quoted here-strings: no $ / ` expansion) ----
$s189 = @'
Followed upPaul, we don't have the time to do scripting, what we are looking for is sample reports like what are available in PVSyst. I agree that it would be useful to provide the report editor available but you will get lots of request for assistance since scripting knowledge is required. 
Also when under the 'Electric Load' tab you get the error below
Could not evaluate callback function:visualize_load_data->on_change
[2]: failed to evaluate function call argument 6 to 'dview()'
Concerning the need for load;
The part that doesn't make sense to me is how to set the parameters for large projects that we have no load, just generation. there is a rate plan we have input and the value of the power changes with the load size. so we have to create a load profile that is large enough to use the generation.
Our buy rate is $0.00 because we are generation only. I have attached the project, If you look at Option A Fixed, that profile has the load set to 0. Please help
'@
$s190 = @'
Emailed to Janine
'@
$s191 = @'
the help data for 'performance model outputs' under "Results" does not work
Also 'Time dependent Pricing Overview' also comes up with a "Page Not Found" error
'@
$s192 = @'
1) In the parametric analysis, the results don't necessarily respond to changes in the main simulation.  In other words, let's say the parametric run tests the effect of different 5 different analysis periods.  I run the main simulation, then run the parametric simulation, no problem.  Now let's say I change the financing terms in the main simulation, and want to re-run the parametric test.  The only way to get it to run again is to change the value in each of the parametric input cells to something different, then change it back.  Alternatively, the number of cells can be decreased (wiping them out), and then increased back.  If I don't do that, re-running the parametric simulation will not produce results that reflect the change in financing that I had entered into the main simulation.
2) The other more minor one that catches me is the inability to specify strings of uneven lengths, when using the advanced residential system design.  Using the "Number of strings in parallel" seems to require the strings to be of the same size.  The alternative using the PV Array Sizing Calculator Algorithm doesn't seem to work well at all, at least in the < 10 kW residential systems I've attempted to model.  It might help if the size field accepted decimal point entry.
'@
$s193 = @'
Forwarded to Aron and Steve to fix parametrics issue
'@
$s194 = @'
Jason Sensibaugh' <sensij@yahoo.com> 
'@

# ===================================================================
# Row 35 (existing row): text updated, status flipped to red, two new
# trailing columns (H/I) added for a follow-up email + date.
# Value is written BEFORE the format-only paste so the copied style
# (date format / wrap / fill / quote-prefix) isn't reset afterwards.
# ===================================================================
$ws.Rows.Item(35).RowHeight = 409.5
$ws.Cells.Item(35,5).Value = $s189

$ws.Cells.Item(2,7).Copy()
$ws.Cells.Item(35,7).PasteSpecial(-4122)

$ws.Cells.Item(35,8).Value = $s190
$ws.Cells.Item(1,4).Copy()
$ws.Cells.Item(35,8).PasteSpecial(-4122)

$ws.Cells.Item(35,9).Value = 41939
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(35,9).PasteSpecial(-4122)

# ===================================================================
# Row 61 (new): Gordon Bleam feedback, marked "Followed up" (green).
# ===================================================================
$ws.Rows.Item(61).RowHeight = 75

$ws.Cells.Item(61,1).Value = 41939
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(61,1).PasteSpecial(-4122)

$ws.Cells.Item(61,2).Value = "Email from SAM"
$ws.Cells.Item(61,3).Value = "Gordon Bleam <gordon@seaportenergy.com>"

$ws.Cells.Item(61,4).Value = $s191
$ws.Cells.Item(1,4).Copy()
$ws.Cells.Item(61,4).PasteSpecial(-4122)

$ws.Cells.Item(61,5).Value = "Followed up"
$ws.Cells.Item(1,4).Copy()
$ws.Cells.Item(61,5).PasteSpecial(-4122)

$ws.Cells.Item(61,6).Value = 41939
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(61,6).PasteSpecial(-4122)

$ws.Cells.Item(3,7).Copy()
$ws.Cells.Item(61,7).PasteSpecial(-4122)

# ===================================================================
# Row 62 (new): Jason Sensibaugh parametric-analysis feedback (red).
# ===================================================================
$ws.Rows.Item(62).RowHeight = 300

$ws.Cells.Item(62,1).Value = 41939
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(62,1).PasteSpecial(-4122)

$ws.Cells.Item(62,2).Value = "Email"

$ws.Cells.Item(62,4).Value = $s192
$ws.Cells.Item(1,4).Copy()
$ws.Cells.Item(62,4).PasteSpecial(-4122)

$ws.Cells.Item(62,5).Value = $s193
$ws.Cells.Item(1,4).Copy()
$ws.Cells.Item(62,5).PasteSpecial(-4122)

$ws.Cells.Item(62,3).Value = $s194
$ws.Cells.Item(17,3).Copy()
$ws.Cells.Item(62,3).PasteSpecial(-4122)

$ws.Cells.Item(62,6).Value = 41939
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(62,6).PasteSpecial(-4122)

$ws.Cells.Item(2,7).Copy()
$ws.Cells.Item(62,7).PasteSpecial(-4122)

# ===================================================================
# Column widths: widen Description-replies column E, add new col I.
# ===================================================================
$ws.Columns.Item(5).ColumnWidth = 44.14
$ws.Columns.Item(9).ColumnWidth = 10.71

# ===================================================================
# View state: keep the frozen pane / selection pinned near the bottom
# of the (now longer) table.
# ===================================================================
$ws.Range("A63").Select()
